$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ruleToJSON")

$ws.Range("U1").Value = "SQL"
$ws.Range("T1").Value = "RULE"

$u2 = '=IF(S2="DISABLE", "", _xlfn.CONCAT(IF(OR(A2=1, A2="1"), "TRUNCATE TABLE rules; ", ""), "INSERT INTO rules (rule_id, rule_msg, correction, pos_example, neg_example) VALUES (", A2, ", ''", SUBSTITUTE(K2, "''", "\''"), "'', ''", SUBSTITUTE(L2, "''", "\''"), "'', ''", SUBSTITUTE(N2, "''", "\''"), "'', ''", SUBSTITUTE(O2, "''", "\''"), "'');"))'
$ws.Range("U2").Formula = $u2

$u3 = '=IF(S3="DISABLE", "", _xlfn.CONCAT(IF(OR(A3=1, A3="1"), "TRUNCATE TABLE rules; ", ""), "INSERT INTO rules (rule_id, rule_msg, correction, pos_example, neg_example) VALUES (", A3, ", ''", SUBSTITUTE(K3, "''", "\''"), "'', ''", SUBSTITUTE(L3, "''", "\''"), "'', ''", SUBSTITUTE(N3, "''", "\''"), "'', ''", SUBSTITUTE(O3, "''", "\''"), "'');"))'
$ws.Range("U3").Formula = $u3

$e8 = '=_xlfn.CONCAT("INSERT INTO rules (rule_id, rule_msg, correction) VALUES (", A2, ", ''", SUBSTITUTE(K2, "''", "\''"), "'', ''", SUBSTITUTE(L2, "''", "\''"), "'', '')")'
$ws.Range("E8").Formula = $e8

# Make the ruleToJSON sheet the active/selected tab, scrolled/selected near
# the newly added SQL column, matching the author's final view state.
$ws.Activate()
$ws.Range("T9").Select()
